$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update columns B,C,D,E,G,H,J,K,L for rows 2-25 (F,I,M,N,O unchanged at 0)
$values = @{
  2 = @{ "B"=1.613131879977203; "C"=0.03254137084699948; "D"=0.5360727841023305; "E"=0.1708840338005366; "G"=4.334273712616181; "H"=3.122628807543634; "J"=0.05983613487714567; "K"=0.9464354327745923; "L"=0.4661599263023533 }
  3 = @{ "B"=1.59662224576212; "C"=0.03108710629339839; "D"=0.5334613920627191; "E"=0.1706915957518476; "G"=4.269737376638943; "H"=3.096292906208873; "J"=0.0598931742378328; "K"=0.9307856508188763; "L"=0.4643624135668034 }
  4 = @{ "B"=1.587566230612936; "C"=0.03017347245391022; "D"=0.5320724504480552; "E"=0.1706252048773784; "G"=4.231067275010247; "H"=3.080752974257052; "J"=0.05993340117343315; "K"=0.9219384899686958; "L"=0.4634856035968582 }
  5 = @{ "B"=1.584147832129872; "C"=0.02979589736936816; "D"=0.5315604184872456; "E"=0.1706111845121221; "G"=4.215548349340452; "H"=3.074578592190903; "J"=0.05995110500261447; "K"=0.9185247973528163; "L"=0.463185395522288 }
  6 = @{ "B"=1.58359664350121; "C"=0.02973288192321633; "D"=0.5314786571502879; "E"=0.1706096442043545; "G"=4.212985882648127; "H"=3.073562892224004; "J"=0.05995412396037558; "K"=0.9179695290710441; "L"=0.4631389963793993 }
  7 = @{ "B"=1.587519027327659; "C"=0.03016840172158197; "D"=0.5320653264105886; "E"=0.1706249629945376; "G"=4.230857012885906; "H"=3.080669063929378; "J"=0.05993363462240175; "K"=0.9218916759184737; "L"=0.4634813236245421 }
  8 = @{ "B"=1.607215099580458; "C"=0.03204420571842448; "D"=0.5351278884634354; "E"=0.1708069468786846; "G"=4.311822766586147; "H"=3.11341717883866; "J"=0.0598547232127169; "K"=0.9408812624154166; "L"=0.4654930852100563 }
  9 = @{ "B"=1.654413203001923; "C"=0.03556087880696168; "D"=0.5428342081960267; "E"=0.1715739679647328; "G"=4.478224038981551; "H"=3.182654660653043; "J"=0.059741178643268; "K"=0.9841683127907004; "L"=0.4712368892178347 }
  10 = @{ "B"=1.694320666551278; "C"=0.03804962614378837; "D"=0.5495327882580199; "E"=0.1723869781512022; "G"=4.605207976676638; "H"=3.236612618881139; "J"=0.05968275524433153; "K"=1.019669660551727; "L"=0.4765528716263105 }
  11 = @{ "B"=1.713613018287987; "C"=0.03916200327206099; "D"=0.5528053819309093; "E"=0.1728109128445965; "G"=4.664021043613502; "H"=3.261837167379724; "J"=0.0596615807053773; "K"=1.036626118356025; "L"=0.4792092389602089 }
  12 = @{ "B"=1.72108220831015; "C"=0.03958045059507498; "D"=0.5540770231017405; "E"=0.1729792124715281; "G"=4.686443790066505; "H"=3.271487077020822; "J"=0.05965433726594505; "K"=1.043163227420166; "L"=0.4802493545733597 }
  13 = @{ "B"=1.719466309863378; "C"=0.03949045346912783; "D"=0.5538017126957016; "E"=0.1729426210095752; "G"=4.68160789745707; "H"=3.269404436295645; "J"=0.05965586283796398; "K"=1.041750180847401; "L"=0.4800238261839098 }
  14 = @{ "B"=1.714224234793477; "C"=0.03919648469101134; "D"=0.5529093518049706; "E"=0.1728246034156946; "G"=4.665862734737175; "H"=3.262629106992449; "J"=0.05966096926431419; "K"=1.037161603720648; "L"=0.4792941245127906 }
  15 = @{ "B"=1.711034615087243; "C"=0.03901605922456497; "D"=0.5523669711550951; "E"=0.1727533250337139; "G"=4.656238126251822; "H"=3.258491783566114; "J"=0.05966419795440636; "K"=1.03436608501309; "L"=0.4788516150222506 }
  16 = @{ "B"=1.693082763035278; "C"=0.03797653676435431; "D"=0.5493234492900001; "E"=0.1723603603227701; "G"=4.601385557823534; "H"=3.234977820925423; "J"=0.05968424759355884; "K"=1.018577754409478; "L"=0.4763840602158211 }
  17 = @{ "B"=1.682361352523174; "C"=0.0373338009414681; "D"=0.5475140510740744; "E"=0.1721331337140484; "G"=4.568004191476035; "H"=3.220726843036402; "J"=0.05969792970527266; "K"=1.009098785297027; "L"=0.4749312553569212 }
  18 = @{ "B"=1.676301822194972; "C"=0.03696225237164441; "D"=0.5464945480974848; "E"=0.1720075311956428; "G"=4.548902648288362; "H"=3.21259393085495; "J"=0.05970630792538856; "K"=1.003722660706245; "L"=0.4741180533200833 }
  19 = @{ "B"=1.674268573704268; "C"=0.03683613055769541; "D"=0.5461530064809494; "E"=0.1719658793554224; "G"=4.542452090457715; "H"=3.209851227469528; "J"=0.05970923206570333; "K"=1.001915436225801; "L"=0.4738465675123678 }
  20 = @{ "B"=1.683491576903236; "C"=0.03740241376473108; "D"=0.5477044691363915; "E"=0.1721567954493288; "G"=4.571547490396; "H"=3.222237271121628; "J"=0.05969642059400293; "K"=1.010099979105632; "L"=0.4750835895015513 }
  21 = @{ "B"=1.715759519992815; "C"=0.03928290559744596; "D"=0.5531705812318819; "E"=0.1728590574177069; "G"=4.670483351474672; "H"=3.26461652541974; "J"=0.05965944836817805; "K"=1.038506229501763; "L"=0.4795075274565903 }
  22 = @{ "B"=1.737802013272244; "C"=0.04049570471396891; "D"=0.5569317308538189; "E"=0.1733632790821744; "G"=4.736027218235336; "H"=3.292884720693678; "J"=0.05963980068304942; "K"=1.057747915998334; "L"=0.4825981977405576 }
  23 = @{ "B"=1.725950296504493; "C"=0.03984987551398689; "D"=0.5549070745142615; "E"=0.1730900302677831; "G"=4.700964079622395; "H"=3.277745111876072; "J"=0.05964987447493897; "K"=1.047416337691942; "L"=0.4809304162840959 }
  24 = @{ "B"=1.682980277413634; "C"=0.03737140023928021; "D"=0.5476183165148569; "E"=0.1721460823008449; "G"=4.569945285149288; "H"=3.221554219185151; "J"=0.05969710126779937; "K"=1.009647110568267; "L"=0.4750146505977852 }
  25 = @{ "B"=1.640726572948353; "C"=0.03462646620668863; "D"=0.5405673861048257; "E"=0.171322617997717; "G"=4.432383360631775; "H"=3.163383758584473; "J"=0.05976749705834017; "K"=0.9718093786131021; "L"=0.4694905344573215 }
}

foreach ($row in $values.Keys) {
  $rowData = $values[$row]
  foreach ($col in $rowData.Keys) {
    $ws.Range("$col$row").Value = $rowData[$col]
  }
}
